$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5: "An Elegant Defense" ---------------------------------------
$ws.Cells.Item(5, 1).Value = "An Elegant Defense"          # A5  (new string idx 14)
$ws.Cells.Item(5, 2).Value = "Matt Richtel"                 # B5  (new string idx 15)

# --- New "Type" column (F) --------------------------------------------------
$ws.Cells.Item(3, 6).Value = "Audio"                         # F3  (new string idx 16)
$ws.Cells.Item(2, 6).Value = "Hard Copy"                     # F2  (new string idx 17)
$ws.Cells.Item(1, 6).Value = "Type"                          # F1  (new string idx 18)

# --- Remaining Tags cell for the new row ------------------------------------
$ws.Cells.Item(5, 5).Value = "immunology;health;medicine"    # E5  (new string idx 19)

# --- New "Length" column (G) ------------------------------------------------
$ws.Cells.Item(1, 7).Value = "Length"                        # G1  (new string idx 20)
$ws.Cells.Item(3, 7).Value = "8 Hrs 36 Mins"                 # G3  (new string idx 21)
$ws.Cells.Item(5, 7).Value = "12 Hrs 38 Mins"                # G5  (new string idx 22)
$ws.Cells.Item(4, 7).Value = "337 Pages"                     # G4  (new string idx 23)
$ws.Cells.Item(2, 7).Value = "304 Pages"                     # G2  (new string idx 24)

# --- Fill remaining reused Type values --------------------------------------
$ws.Cells.Item(5, 6).Value = "Audio"                          # F5 (reuse idx 16)
$ws.Cells.Item(4, 6).Value = "Hard Copy"                      # F4 (reuse idx 17)

# --- Start/Finish dates for the new row, matching existing date styling ----
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)                           # -4122 = xlPasteFormats
$ws.Range("C5").Value = 43832                                 # 1/2/2020

$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = 43840                                 # 1/10/2020
$excel.CutCopyMode = 0

# --- Resize columns to fit the new content ----------------------------------
$ws.Range("A1:G1").EntireColumn.AutoFit() | Out-Null

# --- Restore the editor's last selected cell --------------------------------
$ws.Range("G3").Select() | Out-Null
